# Added test for WhatsApp: append new submission-time rows to each sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Submit orders" -> rows 84, 85 appended (was A1:E83) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A84").Value = "10.13.2022 19:46 (Kyiv+Israel) 16:46 (UTC) 01:46 (Japan) 22:16 (India)"
$ws1.Range("B84").Value = 1.586
$ws1.Range("C84").Value = -0.8590000000000001
$ws1.Range("D84").Value = "***"
$ws1.Range("E84").Value = "***"

$ws1.Range("A85").Value = "10.13.2022 20:06 (Kyiv+Israel) 17:06 (UTC) 02:06 (Japan) 22:36 (India)"
$ws1.Range("B85").Value = 0.802
$ws1.Range("C85").Value = -0.07500000000000007
$ws1.Range("D85").Value = "***"
$ws1.Range("E85").Value = "***"

# --- Sheet 2: "Submit internet survey" -> rows 79, 80 appended (was A1:E78) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A79").Value = "10.13.2022 19:49 (Kyiv+Israel) 16:49 (UTC) 01:49 (Japan) 22:19 (India)"
$ws2.Range("B79").Value = 0.664
$ws2.Range("C79").Value = -0.03300000000000003
$ws2.Range("D79").Value = "***"
$ws2.Range("E79").Value = "***"

$ws2.Range("A80").Value = "10.13.2022 20:09 (Kyiv+Israel) 17:09 (UTC) 02:09 (Japan) 22:39 (India)"
$ws2.Range("B80").Value = 0.881
$ws2.Range("C80").Value = -0.25
$ws2.Range("D80").Value = "***"
$ws2.Range("E80").Value = "***"

# --- Sheet 3: "Submit a phone survey" -> rows 74, 75 appended (was A1:E73) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A74").Value = "10.13.2022 19:51 (Kyiv+Israel) 16:51 (UTC) 01:51 (Japan) 22:21 (India)"
$ws3.Range("B74").Value = 1.356
$ws3.Range("C74").Value = -0.252
$ws3.Range("D74").Value = "***"
$ws3.Range("E74").Value = "***"

$ws3.Range("A75").Value = "10.13.2022 20:11 (Kyiv+Israel) 17:11 (UTC) 02:11 (Japan) 22:41 (India)"
$ws3.Range("B75").Value = 1.406
$ws3.Range("C75").Value = -0.3019999999999998
$ws3.Range("D75").Value = "***"
$ws3.Range("E75").Value = "***"

# --- Sheet 4: "Checkertificate" -> row 90 appended (was A1:E89) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A90").Value = "10.13.2022 20:14 (Kyiv+Israel) 17:14 (UTC) 02:14 (Japan) 22:44 (India)"
$ws4.Range("B90").Value = 0.754
$ws4.Range("C90").Value = -0.08899999999999997
$ws4.Range("D90").Value = "***"
$ws4.Range("E90").Value = "***"
